$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-25 23:48:43'
$ws.Range('E3').Value = '2026-02-25 23:48:45'
$ws.Range('O3').Value = '3.5 °C'
$ws.Range('E4').Value = '2026-02-25 23:48:48'
$ws.Range('J4').Value = '1022.3 hPa'
$ws.Range('E5').Value = '2026-02-25 23:48:51'
$ws.Range('E6').Value = '2026-02-25 23:48:53'
$ws.Range('J6').Value = '1022.2 hPa'
$ws.Range('E7').Value = '2026-02-25 23:48:56'
$ws.Range('E8').Value = '2026-02-25 23:48:59'
$ws.Range('J8').Value = '1021.5 hPa'
$ws.Range('E9').Value = '2026-02-25 23:49:02'
$ws.Range('E10').Value = '2026-02-25 23:49:04'
$ws.Range('E11').Value = '2026-02-25 23:49:07'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '65%'
$ws.Range('O11').Value = '8.4 °C'
$ws.Range('E12').Value = '2026-02-25 23:49:10'
$ws.Range('E13').Value = '2026-02-25 23:49:12'
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H13').Value = '66%'
$ws.Range('J13').Value = '1023.4 hPa'
$ws.Range('O13').Value = '6.4 °C'
$ws.Range('E14').Value = '2026-02-25 23:49:15'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '92%'
$ws.Range('E15').Value = '2026-02-25 23:49:17'
$ws.Range('E16').Value = '2026-02-25 23:49:20'
$ws.Range('N16').Value = '0.0 °C 23:15 TU'
$ws.Range('E17').Value = '2026-02-25 23:49:23'
$ws.Range('O17').Value = '8.8 °C'
$ws.Range('E18').Value = '2026-02-25 23:49:25'
$ws.Range('E19').Value = '2026-02-25 23:49:28'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '52%'
$ws.Range('N19').Value = '7.5 °C 23:08 TU'
$ws.Range('O19').Value = '11.7 °C'
$ws.Range('E20').Value = '2026-02-25 23:49:31'
$ws.Range('N20').Value = '-2.1 °C 23:27 TU'
$ws.Range('O20').Value = '2.3 °C'
$ws.Range('E21').Value = '2026-02-25 23:49:34'
$ws.Range('E22').Value = '2026-02-25 23:49:37'
$ws.Range('O22').Value = '2.2 °C'
$ws.Range('E23').Value = '2026-02-25 23:49:39'
$ws.Range('E24').Value = '2026-02-25 23:49:42'
$ws.Range('J24').Value = '1020.5 hPa'
$ws.Range('E25').Value = '2026-02-25 23:49:45'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '36%'
$ws.Range('O25').Value = '4.9 °C'
$ws.Range('E26').Value = '2026-02-25 23:49:48'
$ws.Range('J26').Value = '1020.2 hPa'
$ws.Range('O26').Value = '9.6 °C'
$ws.Range('E27').Value = '2026-02-25 23:49:51'
$ws.Range('O27').Value = '4.9 °C'
$ws.Range('E28').Value = '2026-02-25 23:49:53'
$ws.Range('E29').Value = '2026-02-25 23:49:56'
$ws.Range('E30').Value = '2026-02-25 23:49:59'
$ws.Range('J30').Value = '1022.3 hPa'
$ws.Range('E31').Value = '2026-02-25 23:50:02'
$ws.Range('J31').Value = '1021.9 hPa'
$ws.Range('E32').Value = '2026-02-25 23:50:05'
$ws.Range('O32').Value = '8.7 °C'
$ws.Range('E33').Value = '2026-02-25 23:50:07'
$ws.Range('J33').Value = '1021.8 hPa'
$ws.Range('O33').Value = '8.2 °C'
$ws.Range('E34').Value = '2026-02-25 23:50:10'
$ws.Range('E35').Value = '2026-02-25 23:50:13'
$ws.Range('J35').Value = '1020.1 hPa'
$ws.Range('N35').Value = '8.1 °C 23:25 TU'
$ws.Range('O35').Value = '12.2 °C'
$ws.Range('E36').Value = '2026-02-25 23:50:16'
$ws.Range('J36').Value = '1022.4 hPa'
$ws.Range('E37').Value = '2026-02-25 23:50:19'
$ws.Range('J37').Value = '1024.0 hPa'
$ws.Range('O37').Value = '6.3 °C'
$ws.Range('E38').Value = '2026-02-25 23:50:21'
$ws.Range('E39').Value = '2026-02-25 23:50:23'
$ws.Range('E40').Value = '2026-02-25 23:50:26'
$ws.Range('H40').NumberFormat = '@'
$ws.Range('H40').Value = '63%'
$ws.Range('J40').Value = '1022.4 hPa'
$ws.Range('O40').Value = '9.1 °C'
$ws.Range('E41').Value = '2026-02-25 23:50:29'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '92%'
$ws.Range('J41').Value = '1021.5 hPa'
$ws.Range('O41').Value = '11.5 °C'
$ws.Range('E42').Value = '2026-02-25 23:50:32'
$ws.Range('O42').Value = '11.5 °C'
$ws.Range('E43').Value = '2026-02-25 23:50:34'
$ws.Range('O43').Value = '9.4 °C'
$ws.Range('E44').Value = '2026-02-25 23:50:37'
$ws.Range('N44').Value = '-1.8 °C 23:09 TU'
$ws.Range('O44').Value = '1.9 °C'
$ws.Range('E45').Value = '2026-02-25 23:50:40'
$ws.Range('J45').Value = '1020.5 hPa'
$ws.Range('O45').Value = '10.5 °C'
$ws.Range('E46').Value = '2026-02-25 23:50:42'
